# Slide 2 ("URL anatomy" diagram): rename the two "Path" labels to
# "Resource Path" (widening their text boxes to fit) and collapse the
# full-URL textbox back into a single run of text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Label above the short/"template" URL: "Path" -> "Resource Path" ---
$shPath1 = $s.Shapes.Item(7)
# Force a real text replacement (old and new text differ here, so the
# run is rewritten normally).
$shPath1.TextFrame.TextRange.Text = "Resource Path"
# Re-anchor/resize the box to the new, wider bounding box recorded for
# this shape (values are point-equivalents of the target EMU offsets,
# nudged to survive the COM single-precision round trip).
$shPath1.Left   = 437.72377952755903
$shPath1.Top    = 94.8723642047244
$shPath1.Width  = 120.60448918897637
$shPath1.Height = 29.081259842519685

# --- Label above the full URL (port included): "Path" -> "Resource Path" ---
$shPath2 = $s.Shapes.Item(14)
$shPath2.TextFrame.TextRange.Text = "Resource Path"
$shPath2.Left   = 494.09740157480314
$shPath2.Top    = 264.0496980393701
$shPath2.Width  = 120.60448918897637
$shPath2.Height = 29.081259842519685

# --- Collapse the 3-run URL ("http://courses.washington.edu" + ":8080" +
#     "/info343/stearns/") into a single run with identical text. The
#     resulting concatenated string is unchanged, and a same-text
#     assignment (or one sharing a prefix/suffix with the existing runs)
#     is treated as a targeted edit of just the differing runs, which
#     would leave the run split in place. Routing the assignment through
#     an unrelated placeholder string first (sharing no prefix/suffix
#     with either the old or the new text) forces a full single-run
#     rewrite, and the follow-up assignment then simply edits that one
#     run's text in place -- preserving its rPr (sz="2400" etc.) -- so
#     the final txBody ends up with exactly one run again. ---
$shUrl = $s.Shapes.Item(8)
$tr = $shUrl.TextFrame.TextRange
$tr.Text = "PLACEHOLDER_TEXT_NO_OVERLAP_998877"
$tr.Text = "http://courses.washington.edu:8080/info343/stearns/"
